$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table on slide 5: switch its table style from the custom "Table_0" style
#    to the built-in "No Style, No Grid" style.
# ---------------------------------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{3BF942C0-C302-43F4-B507-A5C375DDA39D}")

# ---------------------------------------------------------------------------
# 2) Switch the deck's applied design back from the "Integral" (Red Violet)
#    theme colors to the original "Office Theme" colors.
# ---------------------------------------------------------------------------
$officeColors = @(
    0x000000,  # 1  dk1
    0xFFFFFF,  # 2  lt1
    0x6A5444,  # 3  dk2
    0xE6E6E7,  # 4  lt2
    0xD59B5B,  # 5  accent1
    0x317DED,  # 6  accent2
    0xA5A5A5,  # 7  accent3
    0x00C0FF,  # 8  accent4
    0xC47244,  # 9  accent5
    0x47AD70,  # 10 accent6
    0xC16305,  # 11 hlink
    0x724F95   # 12 folHlink
)

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $scheme = $slide.ThemeColorScheme
    for ($c = 1; $c -le $scheme.Count; $c++) {
        $scheme.Colors($c).RGB = $officeColors[$c - 1]
    }
}
